$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "youth basketball leggings with knee pads"
$ws.Range("A2").Value = "youth basketball compression pants with knee pads"
$ws.Range("A3").Value = "softball sliding knee pads"
$ws.Range("A4").Value = "baseball sliding knee pads"
$ws.Range("A5").Value = "hex pads basketball knee youth"
$ws.Range("A6").Value = "youth basketball knee pads"
$ws.Range("A7").Value = "sliding knee pads softball"
$ws.Range("A8").Value = "youth black knee pads for basketball"
$ws.Range("A9").Value = "crossfit leggings"
$ws.Range("A10").Value = "knee pads for basketball"
$ws.Range("A11").Value = "basketball pads"
$ws.Range("A12").Value = "basketball knee"
$ws.Range("A13").Value = "basketball leggings for boys"
$ws.Range("A14").Value = "basketball leg pads tights"
$ws.Range("A15").Value = "compression basketball pants men"
$ws.Range("A16").Value = "crossfit knee sleeves"
$ws.Range("A17").Value = "padded leg sleeves for basketball youth"
$ws.Range("A18").Value = "hex knee pad"
$ws.Range("A19").Value = "compression pants men capri"
$ws.Range("A20").Value = "knee sleeves for crossfit"
$ws.Range("A21").Value = "basketball leg compression pants"
$ws.Range("A22").Value = "compression leggings for men"
$ws.Range("A23").Value = "snowboarding knee pads men"
$ws.Range("A24").Value = "mens tights and leggings"
$ws.Range("A25").Value = "adult football pants with pads"
$ws.Range("A26").Value = "volleyball knee pads xxl"
$ws.Range("A27").Value = "compression basketball pants boys"
$ws.Range("A28").Value = "basketball knee compression sleeve"
$ws.Range("A29").Value = "padded knee compression sleeve"
$ws.Range("A30").Value = "knee pads basketball girls"
$ws.Range("A31").Value = "knee pads for work pants"
$ws.Range("A32").Value = "knee brace wrestling"
$ws.Range("A33").Value = "crossfit knee support"
$ws.Range("A34").Value = "men compression leggings"
$ws.Range("A35").Value = "basketball training gear"
$ws.Range("A36").Value = "sports knee pad"
$ws.Range("A37").Value = "hex padded compression leg sleeve"
$ws.Range("A38").Value = "knee sleeve weightlifting men"
$ws.Range("A39").Value = "hex pads"
$ws.Range("A40").Value = "cycling capris padded"
$ws.Range("A41").Value = "knee pads volleyball men"
$ws.Range("A42").Value = "boy compression pants basketball"
$ws.Range("A43").Value = "crossfit knee compression sleeve"
$ws.Range("A44").Value = "crossfit knee sleeve"
$ws.Range("A45").Value = "mens tights"
$ws.Range("A46").Value = "protector de rodillas basketball"
$ws.Range("A47").Value = "workout pads"
$ws.Range("A48").Value = "mens athletic tights"
$ws.Range("A49").Value = "athletic capri"
$ws.Range("A50").Value = "mens basketball leggings"
$ws.Range("A51").Value = "mens sliding pants baseball"
$ws.Range("A52").Value = "padded compression pants"
$ws.Range("A53").Value = "sliding pants baseball"
$ws.Range("A54").Value = "tights capri men"
$ws.Range("A55").Value = "tights with knee pads basketball"
$ws.Range("A56").Value = "workout knee pads for men"
$ws.Range("A57").Value = "workout sliding pads"
$ws.Range("A58").Value = "wrestling knee pad sleeve"
$ws.Range("A59").Value = "youth knee pads basketball"
$ws.Range("A60").Value = "youth wrestling knee pads"
$ws.Range("A61").Value = "leggings knee pads"
$ws.Range("A62").Value = "xl knee pads"
$ws.Range("A63").Value = "gel knee pads wrestling"
$ws.Range("A64").Value = "men s leggings compression"
$ws.Range("A65").Value = "capri tights for men"
$ws.Range("A66").Value = "padded knee pads"
$ws.Range("A67").Value = "football knee pads youth"
$ws.Range("A68").Value = "knee compression sleeve cycling"
$ws.Range("A69").Value = "volleyball knee pads for girls"
$ws.Range("A70").Value = "youth football leggings"
$ws.Range("A71").Value = "youth basketball leggings boys"
$ws.Range("A72").Value = "wrestling knee sleeve"
$ws.Range("A73").Value = "compression sleeve youth"
$ws.Range("A74").Value = "workout capri pants"
$ws.Range("A75").Value = "workout tights"
$ws.Range("A76").Value = "basketball pants"
$ws.Range("A77").Value = "compression leggings men"
$ws.Range("A78").Value = "youth knee pads"
$ws.Range("A79").Value = "basketball leggins with knee pads"
$ws.Range("A80").Value = "knee pads black mens"
$ws.Range("A81").Value = "basketball knee pads white"
$ws.Range("A82").Value = "knee pads white basketball"
$ws.Range("A83").Value = "black workout capris"
$ws.Range("A84").Value = "black workout leggings capri"
$ws.Range("A85").Value = "workout pads for knees"
$ws.Range("A86").Value = "knee pads for basketball blue"
$ws.Range("A87").Value = "knee tights for men"
$ws.Range("A88").Value = "mcdavid padded leg sleeve"
$ws.Range("A89").Value = "teen leggings"
$ws.Range("A90").Value = "workout pants for men"
$ws.Range("A91").Value = "mens training pants"
$ws.Range("A92").Value = "basketball knee pads women"
$ws.Range("A93").Value = "legging for men workout"
$ws.Range("A94").Value = "mens athletic leggings"
$ws.Range("A95").Value = "mens capri"
$ws.Range("A96").Value = "workout support"
$ws.Range("A97").Value = "athletic legging"
$ws.Range("A98").Value = "basketball knee pads boys youth"
$ws.Range("A99").Value = "basketball spandex knee pads"
$ws.Range("A100").Value = "training gear six pad"
